$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(47)
$pStart = $p.Range.Start
$erRng = $d.Range($pStart + 20, $pStart + 22)
Write-Output ("erRng text=[" + $erRng.Text + "]")
$erRng.Text = [string][char]0xE9
